$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (D) / Volume(1h) (E) readings from this run's scrape.
# Each written cell is forced to Text format first so Excel does not
# auto-coerce a numeric-looking string (e.g. "91.50") into a Number and
# silently drop trailing zeros / switch to scientific notation - the
# source column is plain text, same as the rest of the sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.107.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.968.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -6.03%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.53"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.35%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4985"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4215"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.18"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09041"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.101"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.09"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -6.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.887"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -7.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.442"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.939.36"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001101"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.50"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -9.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06666"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.19"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -9.01%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.944"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.110.94"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.94"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.288"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.62"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.06"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.147"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -11.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.263"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -9.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.25"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.044"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -7.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09843"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.535"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.790"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -7.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.689"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02424"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.293"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.936"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -11.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06311"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6454"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.48"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -8.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1991"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -9.77%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6215"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.44"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.174"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.278"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.470"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000321"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -8.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06897"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.103"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -8.93%  "
